# Update the cryptos list snapshot.
# - Refresh Price (D) / Volume(1h) (E) values for every row.
# - Row 28/29: Toncoin and BitcoinCash swap places (Toncoin moves up to 28).
# - A new coin "Frax" is inserted at row 34, pushing HuobiToken..Aave down
#   by one row and dropping NEARProtocol off the bottom of the list (row 51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price/volume refresh for rows whose coin & link stay put ---
$simple = @{
    2  = @('26.829.37', '  +1.50%  ')
    3  = @('1.726.14',  '  +0.09%  ')
    5  = @('241.20',    '  -0.73%  ')
    6  = @('0.9983',    '  -0.14%  ')
    7  = @('0.4863',    '  -1.06%  ')
    8  = @('0.2588',    '  -1.02%  ')
    9  = @('0.06206',   '  +0.02%  ')
    10 = @('1.722.15',  '  -0.12%  ')
    11 = @('15.97',     '  +3.39%  ')
    12 = @('0.06911',   '  -1.24%  ')
    13 = @('0.6068',    '  +1.26%  ')
    14 = @('4.473',     '  -1.96%  ')
    15 = @('77.09',     '  -0.30%  ')
    16 = @('0.9982',    '  -0.17%  ')
    17 = @('26.607.38', '  +0.65%  ')
    18 = @('0.9983',    '  -0.14%  ')
    19 = @('0.000007163','  -0.13%  ')
    20 = @('11.44',     '  +0.86%  ')
    21 = @('1.954.31',  '  +0.65%  ')
    22 = @('4.421',     '  -1.27%  ')
    23 = @('8.573',     '  -0.23%  ')
    24 = @('5.079',     '  -1.59%  ')
    25 = @('137.70',    '  -0.22%  ')
    26 = @('15.28',     '  +0.24%  ')
    27 = @('1.772',     '  +3.33%  ')
    30 = @('3.931',     '  -0.62%  ')
    31 = @('0.07973',   '  +0.11%  ')
    32 = @('3.679',     '  +0.29%  ')
    33 = @('0.04515',   '  -0.19%  ')
}

foreach ($row in $simple.Keys) {
    $vals = $simple[$row]
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $vals[0]
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $vals[1]
}

# --- Rows whose Coin / Link / Price / Volume all change ---
# column order: B(Coin), C(Link), D(Price), E(Volume)
$full = @{
    28 = @('Toncoin',          'https://coinranking.com/coin/67YlI0K1b+toncoin-ton',                    '1.379',   '  -1.47%  ')
    29 = @('BitcoinCash',      'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch',             '106.06',  '  -0.87%  ')
    34 = @('Frax',             'https://coinranking.com/coin/KfWtaeV1W+frax-frax',                       '0.9977',  '  -0.13%  ')
    35 = @('HuobiToken',       'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht',               '2.600',   '  -0.09%  ')
    36 = @('ARBITRUM',         'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb',                    '1.011',   '  +1.50%  ')
    37 = @('ImmutableX',       'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx',                  '0.6243',  '  -0.33%  ')
    38 = @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt',            '0.9316',  '  +0.43%  ')
    39 = @('RenderToken',      'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr',            '2.033',   '  +3.85%  ')
    40 = @('MXToken',          'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx',                   '2.438',   '  +1.98%  ')
    41 = @('PaxDollar',        'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp',                  '0.9980',  '  -0.16%  ')
    42 = @('VeChain',          'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet',                 '0.01496', '  +0.54%  ')
    43 = @('FraxShare',        'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs',                   '5.644',   '  +5.63%  ')
    44 = @('Quant',            'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt',                   '99.30',   '  -0.36%  ')
    45 = @('TheSandbox',       'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand',                 '0.3834',  '  -0.37%  ')
    46 = @('Aptos',            'https://coinranking.com/coin/HGYj5JCv5+aptos-apt',                       '6.843',   '  +1.29%  ')
    47 = @('Algorand',         'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo',               '0.1160',  '  -0.53%  ')
    48 = @('Cronos',           'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro',                   '0.05401', '  +0.62%  ')
    49 = @('EnergySwap',       'https://coinranking.com/coin/SbWqqTui-+energyswap-ens',                  '7.911',   '  +2.33%  ')
    50 = @('Elrond',           'https://coinranking.com/coin/omwkOTglq+elrond-egld',                     '30.14',   '  +0.04%  ')
    51 = @('Aave',             'https://coinranking.com/coin/ixgUfzmLR+aave-aave',                       '51.57',   '  +1.34%  ')
}

foreach ($row in $full.Keys) {
    $vals = $full[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $dCell = $ws.Cells.Item($row, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $vals[2]
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $vals[3]
}
